$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill missing location for row 12 ---
$ws.Range("G12").Value = "BODEGA QUITO"

# --- Add new row 13 by copying formatting from row 12, then updating values ---
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial()

$ws.Range("A13").Value = "Prueba 2"
$ws.Range("B13").Value = 45454545
$ws.Range("D13").Value = 15245863
$ws.Range("E13").Value = "Puma"
$ws.Range("F13").Value = "8546450004524365875461320323000002450"
$ws.Range("G13").Value = "BODEGA QUITO"

# --- Fix typo in F10 barcode value ---
$ws.Range("F10").Value = "8546454879452365875461320323000002450"

# --- Header renames ---
$ws.Range("G1").Value = "Ubicacion"
$ws.Range("B1").Value = "Codigo"

# --- Update selection to match target state ---
$ws.Range("B1").Select()
